$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column B (shifts B..L -> C..M)
$ws.Columns("B").Insert()

# 2. Row 4: change the label text from "ESTABLECIMEINTO :" to "EMPRESA :"
$ws.Range("A4").Value = "EMPRESA :"

# 3. Widen the merged input box for EMPRESA to B4:E4 (absorb newly inserted column)
$ws.Range("B4:E4").Merge()

# 4. Add new "ESTABLECIMEINTO :" label at F4 (re-using the style of A4)
$ws.Range("A4").Copy()
$ws.Range("F4").PasteSpecial(-4122)
$ws.Range("F4").Value = "ESTABLECIMEINTO :"

# 5. Add new merged input box G4:I4 for ESTABLECIMIENTO, matching style of B4
$ws.Range("B4").Copy()
$ws.Range("G4:I4").PasteSpecial(-4122)
$ws.Range("G4:I4").Merge()

# 6. Row 10 header: add new "SUCURSAL" column header at B10
$ws.Range("B10").Value = "SUCURSAL"

# 7. Update selection to B10 (matches target sheetView selection)
$ws.Range("B10").Select()
